# Update "想去人数" (column F) figures across the four sheets to the
# values captured in the latest gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5720
$ws.Range("F6").Value = 1603
$ws.Range("F7").Value = 168
$ws.Range("F8").Value = 3239
$ws.Range("F12").Value = 4518
$ws.Range("F13").Value = 1078
$ws.Range("F15").Value = 1703
$ws.Range("F19").Value = 50
$ws.Range("F20").Value = 172
$ws.Range("F22").Value = 1018
$ws.Range("F23").Value = 304
$ws.Range("F29").Value = 1116
$ws.Range("F31").Value = 83
$ws.Range("F33").Value = 373
$ws.Range("F34").Value = 502
$ws.Range("F36").Value = 1736
$ws.Range("F37").Value = 2246
$ws.Range("F40").Value = 271
$ws.Range("F41").Value = 632
$ws.Range("F42").Value = 364
$ws.Range("F43").Value = 34
$ws.Range("F44").Value = 669
$ws.Range("F45").Value = 29
$ws.Range("F47").Value = 396
$ws.Range("F49").Value = 146

# 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F21").Value = 11

# 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 778

# 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 778
$ws.Range("F5").Value = 5720
$ws.Range("F6").Value = 1603
$ws.Range("F7").Value = 168
$ws.Range("F9").Value = 3239
$ws.Range("F11").Value = 4518
$ws.Range("F12").Value = 1078
$ws.Range("F19").Value = 50
$ws.Range("F20").Value = 172
$ws.Range("F23").Value = 1018
$ws.Range("F24").Value = 304
$ws.Range("F30").Value = 1116
$ws.Range("F33").Value = 503
$ws.Range("F35").Value = 1736
$ws.Range("F41").Value = 271
$ws.Range("F42").Value = 364
$ws.Range("F43").Value = 669
$ws.Range("F45").Value = 396
$ws.Range("F48").Value = 146
